# Fill in the new NXOS equipment row (row 3) on Sheet1, matching the pattern
# already used for row 2 (hostname/ip/os/protocol/username/password/
# enable_password/platform), and move the active selection to H3 - this is
# what happens when a user types the new row in and then tabs/clicks onto H3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "NXOS"
$ws.Range("B3").Value = "192.168.0.212"
$ws.Range("C3").Value = "nxos"
$ws.Range("D3").Value = "ssh"
$ws.Range("E3").Value = "alex"

# Columns F and G follow row 2's look: general (left/default) alignment
# rather than the centered alignment the blank template row started with.
$ws.Range("F3").ClearFormats()
$ws.Range("F3").Value = "alex"

$ws.Range("G3").ClearFormats()
$ws.Range("G3").Value = "alex"

$ws.Range("H3").Value = "nx"

# Matches the saved cursor position (H3) recorded in the worksheet view.
$ws.Range("H3").Select()
